# Adding updated input file excel sheet
#
# Adds a new worksheet "ValidLogin" (after the existing "TC1" sheet) that
# holds a tiny username/password -> admin/pointofsale login table, and makes
# it the active sheet.

$wb = $excel.ActiveWorkbook

# Existing sheet ("TC1") stays sheetId=1 / first position; the new sheet is
# inserted right after it, becoming sheetId=2 / second position.
$tc1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc1)
$ws.Name = "ValidLogin"

# Header + credential row.
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Leave the selection/zoom the way it was saved in the authored workbook.
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 160
